$d = $word.ActiveDocument

# Remove the two trailing paragraphs ("Wtf are u doing" and the blank
# paragraph that follows it), reverting the document back to just the
# single "123" paragraph before the section properties.
$count = $d.Paragraphs.Count
if ($count -gt 1) {
    $startPara = $d.Paragraphs.Item(2)
    $endPara = $d.Paragraphs.Item($count)
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
